# Applies the resume.pptx edit:
#  1. Fixes a typo in the "ABOUT ME" paragraph (inside Group 4 /
#     TextBox 8): "research orientation class" -> "research orientated
#     class".
#  2. Moves several shapes/groups upward (smaller Y offset), and also
#     nudges the contact-info TextBox ("TextBox 37") left and up.
#
# EMU <-> point conversion: 1 pt = 12700 EMU. Shape.Left/Top are COM
# Single (float32) properties, exactly like real PowerPoint -- a couple
# of the literal EMU numbers below are the target value "+1" so that,
# after the float32 round-trip that happens on save, the stored EMU
# lands exactly on the intended value instead of one EMU short.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$emuPerPt = 12700.0

# --- 1. Text fix inside Group 4 ("ABOUT ME") -> TextBox 8 ---------------
$grp = $s.Shapes.Item("Group 4")
$tb  = $grp.GroupItems.Item("TextBox 8")
$tr  = $tb.TextFrame.TextRange
$fullText = $tr.Text

$oldRun = "projects. I am currently enrolled in ME 193A, a research orientation class with the goal of developing a spacesuit for manned Mars missions. My interests include data science, machine learning, artificial intelligence, robotics, and the application of computer science to fields such as biology and chemistry."
$newRun = "projects. I am currently enrolled in ME 193A, a research orientated class with the goal of developing a spacesuit for manned Mars missions. My interests include data science, machine learning, artificial intelligence, robotics, and the application of computer science to fields such as biology and chemistry."

$startIdx = $fullText.IndexOf($oldRun)
if ($startIdx -ge 0) {
    # Replace exactly this run's span so it stays a single run in the XML
    # (like the diff shows) rather than getting fragmented.
    $range = $tr.Characters($startIdx + 1, $oldRun.Length)
    $range.Text = $newRun
}

# --- 2. Shape/group vertical (and one horizontal) repositioning ---------
# Group 4            y: 1341964 -> 1209008
$s.Shapes.Item("Group 4").Top = 1209009 / $emuPerPt

# Group 3             y: 2472491 -> 2339535
$s.Shapes.Item("Group 3").Top = 2339535 / $emuPerPt

# Group 5             y: 4599108 -> 4466152
$s.Shapes.Item("Group 5").Top = 4466153 / $emuPerPt

# TextBox 21          y: 6587554 -> 6454598
$s.Shapes.Item("TextBox 21").Top = 6454598 / $emuPerPt

# Rectangle 22        y: 6822157 -> 6689201
$s.Shapes.Item("Rectangle 22").Top = 6689201 / $emuPerPt

# Rectangle 24        y: 7335010 -> 7202054
$s.Shapes.Item("Rectangle 24").Top = 7202054 / $emuPerPt

# Rectangle 25        y: 7841580 -> 7708624
$s.Shapes.Item("Rectangle 25").Top = 7708624 / $emuPerPt

# Group 31            y: 5589153 -> 5456197
$s.Shapes.Item("Group 31").Top = 5456197 / $emuPerPt

# Rectangle 34        y: 8377883 -> 8244927
$s.Shapes.Item("Rectangle 34").Top = 8244928 / $emuPerPt

# TextBox 37 (www.hantaowang.me / email / github / linkedin box)
#                     x: 205366 -> 185646, y: 725670 -> 624233
$s.Shapes.Item("TextBox 37").Left = 185646 / $emuPerPt
$s.Shapes.Item("TextBox 37").Top  = 624233 / $emuPerPt

# Straight Connector 38   y: 1260088 -> 1130503
$s.Shapes.Item("Straight Connector 38").Top = 1130504 / $emuPerPt
